$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllEntries")

# Insert a blank row below row 17 (the UQAM-CRCM5 entry), shifting rows down,
# then stamp it with row 17's formatting (borders/fonts/number formats) and
# fill in the values for the new CRCM5-SN entry.
$ws.Rows.Item(17).Copy() | Out-Null
$ws.Rows.Item(18).Insert(-4121)  # xlShiftDown
$excel.CutCopyMode = 0

$ws.Range("A17:G17").Copy() | Out-Null
$ws.Range("A18:G18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "CRCM5-SN"
$ws.Range("B18").Value = "UQAM"
$ws.Range("C18").Formula = '=CONCATENATE(B18,"-",A18)'
$ws.Range("D18").Value = "Winger.Katja@uqam.ca"
$ws.Range("E18").Value = "Universite du Quebec a Montreal"
$ws.Range("F18").ClearContents() | Out-Null
$ws.Range("G18").Value = "unrestricted"
$ws.Rows.Item(18).RowHeight = 22.5
